$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Data")

# Clear the contents of the last data row (row 25) while keeping the
# cell formatting (styles) in place - matches Excel's "Clear Contents"
# applied to the row.
$ws.Range("A25:H25").ClearContents()

# Leave the same selection Excel shows after clearing an entire row:
# the whole row selected with the active cell at A25.
$ws.Rows("25:25").Select()
